$d = $word.ActiveDocument

# -----------------------------------------------------------------
# 1) document.xml: the "Heading 2" sample paragraph (bookmark
#    "heading-2") loses its direct rFonts/color run-formatting
#    overrides (falls back to the Heading2 style's own formatting).
# -----------------------------------------------------------------
$h2Para = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -eq " Heading 2 ") {
        $h2Para = $p
        break
    }
}

if ($h2Para -ne $null) {
    $rng = $h2Para.Range
    $fragment = "<?xml version='1.0' standalone='yes'?>" +
        "<?mso-application progid='Word.Document'?>" +
        "<pkg:package xmlns:pkg='http://schemas.microsoft.com/office/2006/xmlPackage'>" +
        "<pkg:part pkg:name='/word/document.xml' pkg:contentType='application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml'>" +
        "<pkg:xmlData>" +
        "<w:document xmlns:w='http://schemas.openxmlformats.org/wordprocessingml/2006/main'>" +
        "<w:body>" +
        "<w:p>" +
        "<w:pPr><w:pStyle w:val='Heading2'/><w:rPr/></w:pPr>" +
        "<w:bookmarkStart w:id='1' w:name='heading-2'/>" +
        "<w:r><w:rPr/><w:t xml:space='preserve'> </w:t></w:r>" +
        "<w:r><w:rPr/><w:t xml:space='preserve'>Heading 2 </w:t></w:r>" +
        "<w:bookmarkEnd w:id='1'/>" +
        "</w:p>" +
        "</w:body></w:document>" +
        "</pkg:xmlData></pkg:part></pkg:package>"
    $rng.InsertXML($fragment)
}

# -----------------------------------------------------------------
# 2) styles.xml: Heading1 style - bottom spacing 0 -> 144 twips (7.2pt)
# -----------------------------------------------------------------
$h1 = $d.Styles("Heading 1")
$h1.ParagraphFormat.SpaceAfter = 7.2

# -----------------------------------------------------------------
# 3) styles.xml: Heading2 style - spacing, alignment, font, color
# -----------------------------------------------------------------
$h2 = $d.Styles("Heading 2")
$h2.ParagraphFormat.SpaceBefore = 0
$h2.ParagraphFormat.SpaceAfter = 36
$h2.ParagraphFormat.Alignment = 1
$h2.Font.Name = "Times New Roman"
$h2.Font.Color = 0

# -----------------------------------------------------------------
# 4) styles.xml: FirstParagraph style - top spacing 0 -> 720 twips (36pt)
# -----------------------------------------------------------------
$fp = $d.Styles("First Paragraph")
$fp.ParagraphFormat.SpaceBefore = 36

Write-Host "done"
